$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-05 Saturday" "2025-04-06 Sunday"

Replace-Text "753×8=" "882×2="
Replace-Text "783×2=" "691×5="
Replace-Text "141×6=" "411×5="
Replace-Text "424×9=" "496×2="
Replace-Text "727×6=" "542×6="

Replace-Text "216×3=" "869×6="
Replace-Text "152×3=" "595×3="
Replace-Text "344×6=" "878×6="
Replace-Text "357×4=" "421×5="
Replace-Text "306×6=" "663×9="

Replace-Text "477×2=" "465×4="
Replace-Text "331×6=" "637×5="
Replace-Text "891×7=" "347×4="
Replace-Text "771×4=" "274×6="
Replace-Text "631×9=" "350×5="

Replace-Text "303×4=" "435×3="
Replace-Text "861×2=" "390×4="
Replace-Text "800×4=" "200×4="
Replace-Text "341×2=" "949×3="
Replace-Text "204×9=" "633×4="

Replace-Text "786×7=" "745×8="
Replace-Text "264×8=" "764×7="
Replace-Text "762×9=" "218×6="
Replace-Text "184×7=" "519×6="
Replace-Text "339×6=" "569×8="
